$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 'Desenvolver nos alunos a capacidade de realizarem práticas rotineiras de laboratório associadas ao desenvolvimento de seu pensamento científico, resolvendo problemas teóricos e práticos, utilizando corretamente os diversos materiais de laboratório e manipulando reagentes químicos com segurança. Dessa forma, os alunos adquirirão experiência  nas várias áreas da química aplicando os conceitos pertinentes'
$ws.Range("C10").Value = 'Desenvolver nos alunos a capacidade de realizarem práticas rotineiras de laboratório associadas ao desenvolvimento de seu pensamento científico, resolvendo problemas teóricos e práticos, utilizando corretamente os diversos materiais de laboratório e manipulando reagentes químicos com segurança. Dessa forma, os alunos adquirirão experiência  nas várias áreas da química aplicando os conceitos pertinentes'

$ws.Range("A13").Clear()
$ws.Range("B13").Value = '5817330 - Larissa de Freitas'
$ws.Range("C13").Value = '5817330 - Larissa de Freitas'
$ws.Rows.Item(13).AutoFit()

$ws.Range("A14").Clear()
$ws.Range("B14").Value = '6310296 - Patrícia Caroline Molgero Da Rós'
$ws.Range("C14").Value = '6310296 - Patrícia Caroline Molgero Da Rós'
$ws.Rows.Item(14).AutoFit()

$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = '1 - Introdução ao Laboratório Químico; 2 - Pesos e Medidas; 3 - Técnicas de Separação de Misturas; 4 - Fenômenos físicos; 5 – Miscibilidade; 6 – Reações Químicas; 7 – Soluções; 8– Titrimetria; 9 – Equilíbrio Químico.'
$ws.Range("C15").Value = '1 - Introdução ao Laboratório Químico; 2 - Pesos e Medidas; 3 - Técnicas de Separação de Misturas; 4 - Fenômenos físicos; 5 – Miscibilidade; 6 – Reações Químicas; 7 – Soluções; 8– Titrimetria; 9 – Equilíbrio Químico.'
$ws.Rows.Item(15).RowHeight = 60

$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = '1 - Introduction to the Chemistry Laboratory; 2 - Weights and measures; 3 - Methods for separating mixtures; 4 - Physical phenomena;  5 - Miscibility; 6 - Chemical Reactions; 7 -  Solutions; 8- Titrimetry; 9 - Chemical Equilibrium.'
$ws.Range("C16").Value = '1 - Introduction to the Chemistry Laboratory; 2 - Weights and measures; 3 - Methods for separating mixtures; 4 - Physical phenomena;  5 - Miscibility; 6 - Chemical Reactions; 7 -  Solutions; 8- Titrimetry; 9 - Chemical Equilibrium.'
$ws.Rows.Item(16).RowHeight = 60

$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = '1 - Introdução ao Laboratório: Noções Elementares de Segurança; Equipamentos Básicos de Laboratório; Equipamentos de Proteção Individual.2 - Pesos e medidas (Tratamento de dados experimentais): Cuidados Gerais com Balanças; Técnicas de Determinação de massa; Exatidão e precisão; Unidades; Algarismos Significativos; Propagação de Erros.3 - Técnicas de Separação de Misturas: Filtração simples; Filtração a vácuo e Decantação.4 - Fenômenos físicos: Construção do Diagrama da mudança do estado físico da água.5 - Miscibilidade e solubilidade: Influência das forças intermoleculares na miscibilidade de líquidos.6 - Reações químicas: Aspectos qualitativos.7 - Soluções: Preparo e padronização de soluções.8 - Titrimetria: Realização de Titulações Ácido-Base; Retrotitulação.9 - Equilíbrio Químico - Preparo de Solução Tampão.'
$ws.Range("C17").Value = '1 - Introdução ao Laboratório: Noções Elementares de Segurança; Equipamentos Básicos de Laboratório; Equipamentos de Proteção Individual.2 - Pesos e medidas (Tratamento de dados experimentais): Cuidados Gerais com Balanças; Técnicas de Determinação de massa; Exatidão e precisão; Unidades; Algarismos Significativos; Propagação de Erros.3 - Técnicas de Separação de Misturas: Filtração simples; Filtração a vácuo e Decantação.4 - Fenômenos físicos: Construção do Diagrama da mudança do estado físico da água.5 - Miscibilidade e solubilidade: Influência das forças intermoleculares na miscibilidade de líquidos.6 - Reações químicas: Aspectos qualitativos.7 - Soluções: Preparo e padronização de soluções.8 - Titrimetria: Realização de Titulações Ácido-Base; Retrotitulação.9 - Equilíbrio Químico - Preparo de Solução Tampão.'
$ws.Rows.Item(17).RowHeight = 120

$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = '1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment.2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation.3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation.4 - Physical phenomena: Water state changes.5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. 6 - Chemical reactions: Qualitative aspects.7 - Solutions: Preparation and standardization of solutions.8 - Titrimetry: Acid-Base Titrations and return-titration.9 - Chemical equilibrium: Buffer solution.'
$ws.Range("C18").Value = '1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment.2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation.3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation.4 - Physical phenomena: Water state changes.5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. 6 - Chemical reactions: Qualitative aspects.7 - Solutions: Preparation and standardization of solutions.8 - Titrimetry: Acid-Base Titrations and return-titration.9 - Chemical equilibrium: Buffer solution.'
$ws.Rows.Item(18).RowHeight = 120

$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()
$ws.Rows.Item(19).AutoFit()

$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Os instrumentos de avaliação utilizados serão duas provas (P1 e P2) e a média dos relatórios (MR). O professor poderá a seu critério utilizar de trabalhos e/ou testes para complementar o método avaliativo.'
$ws.Range("C20").Value = 'Os instrumentos de avaliação utilizados serão duas provas (P1 e P2) e a média dos relatórios (MR). O professor poderá a seu critério utilizar de trabalhos e/ou testes para complementar o método avaliativo.'

$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'A nota final será calculada da seguinte forma: NF = (3xMR + 7xMP)/10 onde NF é a nota final , MR é a média dos relatórios e MP é a média simples das provas.'
$ws.Range("C21").Value = 'A nota final será calculada da seguinte forma: NF = (3xMR + 7xMP)/10 onde NF é a nota final , MR é a média dos relatórios e MP é a média simples das provas.'
$ws.Rows.Item(21).RowHeight = 60

$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'A recuperação será feita por meio de uma prova (PR) para alunos que tenham NF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a nota final (NF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0'
$ws.Range("C22").Value = 'A recuperação será feita por meio de uma prova (PR) para alunos que tenham NF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a nota final (NF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0'
$ws.Rows.Item(22).RowHeight = 60

$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23").Value = 'ASSUMPÇÃO, R. M. V. ; MORITA, T. Manual de soluções reagentes e solventes: padronização, preparação, purificação. São Paulo: Editora Edgard Blucher, 1972.BACCAN, N.; ANDRADE, J. C. O. ; GODINHO, E. S.; BARONE, J. S. Química analítica quantitativa elementar. 2.ed. São Paulo: Edgard Blucher, 1995.BRADY, J; HUMISTON, G. E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1986.BROWN, T. E et al. Química a Ciência Central. 9 ed. São Paulo. Pearson Prentice Hall, 2005-2007.CONSTANTINO, M.G; SILVA, G. V. J. da; DONATE P. M. Fundamentos de química experimental, São Paulo : EDUSP, 2004.MAHAN, B. M.; MYERS, R. J. Química um curso universitário. São Paulo: Ed. Edgard Blucher Ltda, 1993.SILVA, R. R.; BOCCHI, N. ; ROCHA FILHO, R. C. Introdução a química experimental. São Paulo: McGraw-Hill, 1990.'
$ws.Range("C23").Value = 'ASSUMPÇÃO, R. M. V. ; MORITA, T. Manual de soluções reagentes e solventes: padronização, preparação, purificação. São Paulo: Editora Edgard Blucher, 1972.BACCAN, N.; ANDRADE, J. C. O. ; GODINHO, E. S.; BARONE, J. S. Química analítica quantitativa elementar. 2.ed. São Paulo: Edgard Blucher, 1995.BRADY, J; HUMISTON, G. E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1986.BROWN, T. E et al. Química a Ciência Central. 9 ed. São Paulo. Pearson Prentice Hall, 2005-2007.CONSTANTINO, M.G; SILVA, G. V. J. da; DONATE P. M. Fundamentos de química experimental, São Paulo : EDUSP, 2004.MAHAN, B. M.; MYERS, R. J. Química um curso universitário. São Paulo: Ed. Edgard Blucher Ltda, 1993.SILVA, R. R.; BOCCHI, N. ; ROCHA FILHO, R. C. Introdução a química experimental. São Paulo: McGraw-Hill, 1990.'
$ws.Rows.Item(23).RowHeight = 120

$ws.Range("A24").Value = 'Requisitos:'

$ws.Range("B25").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("C25").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Rows.Item(25).RowHeight = 30
